# Updates the Price (D) and Volume(1h) (E) columns for the crypto list, and
# swaps the InternetComputer(DFINITY)/Hedera rows (39<->40) with their new
# data, per the GitHub Actions refresh commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> new D (Price) / E (Volume) values. DNumeric flags rows whose
# new Price text parses as a plain number (single decimal point) so Excel's
# automatic type inference would otherwise silently convert the cell from
# text to a numeric value; for those we force text entry with a leading
# apostrophe so the stored cell keeps its original text semantics.
$changes = @{
    2 = @{ D = '24.382.40'; DNumeric = $false; E = '  +9.01%  ' }
    3 = @{ D = '1.680.58'; DNumeric = $false; E = '  +4.54%  ' }
    4 = @{ D = $null; DNumeric = $false; E = '  -0.50%  ' }
    5 = @{ D = '307.69'; DNumeric = $true; E = '  +5.74%  ' }
    6 = @{ D = '0.9964'; DNumeric = $true; E = '  +0.26%  ' }
    7 = @{ D = '0.3702'; DNumeric = $true; E = '  -0.49%  ' }
    8 = @{ D = '0.3445'; DNumeric = $true; E = '  +1.89%  ' }
    9 = @{ D = '47.93'; DNumeric = $true; E = '  +12.23%  ' }
    10 = @{ D = $null; DNumeric = $false; E = '  +3.06%  ' }
    11 = @{ D = '0.07259'; DNumeric = $true; E = '  +2.27%  ' }
    12 = @{ D = '0.9988'; DNumeric = $true; E = '  -0.27%  ' }
    13 = @{ D = '20.40'; DNumeric = $true; E = '  +2.25%  ' }
    14 = @{ D = '6.108'; DNumeric = $true; E = '  +2.60%  ' }
    15 = @{ D = '6.748'; DNumeric = $true; E = '  +0.98%  ' }
    16 = @{ D = '1.676.87'; DNumeric = $false; E = '  +4.53%  ' }
    17 = @{ D = '0.00001113'; DNumeric = $true; E = '  +2.26%  ' }
    18 = @{ D = '0.9967'; DNumeric = $true; E = '  +0.31%  ' }
    19 = @{ D = '0.06719'; DNumeric = $true; E = '  +1.07%  ' }
    20 = @{ D = '81.34'; DNumeric = $true; E = '  +3.46%  ' }
    21 = @{ D = '16.48'; DNumeric = $true; E = '  +1.45%  ' }
    22 = @{ D = '6.088'; DNumeric = $true; E = '  +0.75%  ' }
    23 = @{ D = '11.96'; DNumeric = $true; E = '  +1.13%  ' }
    24 = @{ D = '24.341.53'; DNumeric = $false; E = '  +8.66%  ' }
    25 = @{ D = '2.437'; DNumeric = $true; E = '  +1.39%  ' }
    26 = @{ D = '3.362'; DNumeric = $true; E = '  -11.44%  ' }
    27 = @{ D = '2.665'; DNumeric = $true; E = '  +6.09%  ' }
    28 = @{ D = '151.91'; DNumeric = $true; E = '  +0.73%  ' }
    29 = @{ D = '19.57'; DNumeric = $true; E = '  +0.06%  ' }
    30 = @{ D = '1.862.14'; DNumeric = $false; E = '  +4.35%  ' }
    31 = @{ D = '127.35'; DNumeric = $true; E = '  +4.87%  ' }
    32 = @{ D = '6.283'; DNumeric = $true; E = '  +4.49%  ' }
    33 = @{ D = '4.037'; DNumeric = $true; E = '  -4.87%  ' }
    34 = @{ D = '0.9725'; DNumeric = $true; E = '  +1.59%  ' }
    35 = @{ D = '1.758'; DNumeric = $true; E = '  +7.28%  ' }
    36 = @{ D = '0.08452'; DNumeric = $true; E = '  +2.22%  ' }
    37 = @{ D = '9.033'; DNumeric = $true; E = '  +3.72%  ' }
    38 = @{ D = '12.27'; DNumeric = $true; E = '  +3.54%  ' }
    41 = @{ D = '0.02330'; DNumeric = $true; E = '  +5.38%  ' }
    42 = @{ D = '1.261'; DNumeric = $true; E = '  +2.23%  ' }
    43 = @{ D = '0.2110'; DNumeric = $true; E = '  +3.91%  ' }
    44 = @{ D = '0.6165'; DNumeric = $true; E = '  +2.55%  ' }
    45 = @{ D = '0.9959'; DNumeric = $true; E = '  +0.40%  ' }
    46 = @{ D = '3.780'; DNumeric = $true; E = '  +3.09%  ' }
    47 = @{ D = '13.02'; DNumeric = $true; E = '  -1.76%  ' }
    48 = @{ D = '0.5936'; DNumeric = $true; E = '  +3.03%  ' }
    49 = @{ D = '127.13'; DNumeric = $true; E = '  +1.26%  ' }
    50 = @{ D = '2.023'; DNumeric = $true; E = '  +2.04%  ' }
    51 = @{ D = '0.07219'; DNumeric = $true; E = '  +4.95%  ' }
}

foreach ($row in $changes.Keys) {
    $info = $changes[$row]
    if ($null -ne $info.D) {
        if ($info.DNumeric) {
            # Leading apostrophe forces Excel to keep the numeric-looking
            # text as a text value instead of silently converting it to a
            # number.
            $ws.Cells.Item($row, 4).Value = "'" + $info.D
        } else {
            $ws.Cells.Item($row, 4).Value = $info.D
        }
    }
    $ws.Cells.Item($row, 5).Value = $info.E
}

# Rows 39 and 40 swap places: InternetComputer(DFINITY) moves from row 39 to
# row 40, and Hedera moves from row 40 to row 39 - each with refreshed
# Price/Volume figures.
$ws.Cells.Item(39, 2).Value = "Hedera"
$ws.Cells.Item(39, 3).Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Cells.Item(39, 4).Value = "'0.06397"
$ws.Cells.Item(39, 5).Value = "  +3.17%  "

$ws.Cells.Item(40, 2).Value = "InternetComputer(DFINITY)"
$ws.Cells.Item(40, 3).Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Cells.Item(40, 4).Value = "'5.343"
$ws.Cells.Item(40, 5).Value = "  -0.05%  "
